$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data

$ws.Range("D2").Value = '66.964.52'
$ws.Range("E2").Value = '  -3.69%  '
$ws.Range("D3").Value = '3.674.64'
$ws.Range("E3").Value = '  -2.94%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.28'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.48'
$ws.Range("E6").Value = '  -6.12%  '
$ws.Range("D7").Value = '3.674.53'
$ws.Range("E7").Value = '  -2.99%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.519'
$ws.Range("E9").Value = '  -1.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.156'
$ws.Range("E10").Value = '  -5.88%  '
$ws.Range("E11").Value = '  -4.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("E12").Value = '  -5.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.46'
$ws.Range("E13").Value = '  -5.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000238'
$ws.Range("E14").Value = '  -6.45%  '
$ws.Range("D15").Value = '4.291.56'
$ws.Range("E15").Value = '  -2.89%  '
$ws.Range("D16").Value = '3.685.21'
$ws.Range("E16").Value = '  -2.78%  '
$ws.Range("D17").Value = '67.018.82'
$ws.Range("E17").Value = '  -3.79%  '
$ws.Range("E18").Value = '  -4.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.05'
$ws.Range("E19").Value = '  -6.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.82'
$ws.Range("E20").Value = '  +1.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '482.61'
$ws.Range("E21").Value = '  -5.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.06'
$ws.Range("E22").Value = '  -5.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.714'
$ws.Range("E23").Value = '  -2.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.36'
$ws.Range("E24").Value = '  -2.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.28'
$ws.Range("E25").Value = '  -7.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000139'
$ws.Range("E26").Value = '  -1.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.04'
$ws.Range("E27").Value = '  -6.37%  '
$ws.Range("E28").Value = '  -0.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.88'
$ws.Range("E29").Value = '  -6.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.90'
$ws.Range("E30").Value = '  -3.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.34'
$ws.Range("E31").Value = '  -7.16%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.58'
$ws.Range("E32").Value = '  +0.94%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.63'
$ws.Range("E33").Value = '  -5.98%  '
$ws.Range("D34").Value = '3.818.98'
$ws.Range("E34").Value = '  -2.81%  '
$ws.Range("D35").Value = '3.617.75'
$ws.Range("E35").Value = '  -2.71%  '
$ws.Range("E36").Value = '  -7.64%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.988'
$ws.Range("E38").Value = '  -5.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.70'
$ws.Range("E39").Value = '  -6.98%  '
$ws.Range("E40").Value = '  -7.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.319'
$ws.Range("E41").Value = '  -6.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '434.84'
$ws.Range("E42").Value = '  -9.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.54'
$ws.Range("E43").Value = '  -2.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.91'
$ws.Range("E44").Value = '  -7.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.74'
$ws.Range("E45").Value = '  -7.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.23'
$ws.Range("E46").Value = '  -4.00%  '
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '141.28'
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("B49").Value = 'Arweave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '39.30'
$ws.Range("E49").Value = '  -10.74%  '
$ws.Range("D50").Value = '2.753.52'
$ws.Range("E50").Value = '  -6.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0343'
$ws.Range("E51").Value = '  -5.64%  '
